$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Merge the two runs
#       "Validate entry by checking for correct data type"
#       " by the use of combo boxes in the forms"
#    into the single new sentence
#       "Use combo boxes when applicable, to help minimize input errors "
#    A temporary trailing sentinel character ("X") is appended to the
#    replacement text so that the collapsed range used to place the
#    bookmark in step 2 does not land exactly on the paragraph-mark
#    boundary (the bookmark API mishandles a collapsed range whose
#    position is precisely "end of paragraph text, just before the
#    paragraph mark"). The sentinel is removed again in step 3.
# ------------------------------------------------------------------
$r = $d.Content
[void]$r.Find.Execute("Validate entry by checking for correct data type by the use of combo boxes in the forms", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Use combo boxes when applicable, to help minimize input errors X", 2)

# ------------------------------------------------------------------
# 2. Re-anchor the existing "_GoBack" bookmark (previously sitting
#    just before the "Output Controls" paragraph) to the end of the
#    text we just inserted, i.e. right after "...input errors " and
#    before the paragraph mark -- where Word leaves it after editing.
#    Bookmarks.Add with an existing bookmark name moves that bookmark.
# ------------------------------------------------------------------
$bmPos = $r.End - 1
$target = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $target)

# ------------------------------------------------------------------
# 3. Drop the sentinel character now that the bookmark is anchored
#    right before it. Deleting text after a bookmark's end does not
#    move the bookmark.
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$d.Range($bm.End, $bm.End + 1).Delete()
